# daily auto push: 2026-02-28 13:45 UTC
# A new observation row (2026/02/28, Sat, hour=19, value=201) was recorded
# and inserted into the daily log, ahead of the existing 2026/12/29 block.
# That pushes every following row down by one (old row 904 -> new row 905,
# ..., old row 945 -> new row 946), and the sheet's used range grows by a
# row (A1:D945 -> A1:D946).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 904 (shifts 904..945 down to 905..946).
$ws.Rows.Item(904).Insert()

# Column A holds the date as literal text (e.g. "2026/12/29"), not a real
# Excel date serial - force the new cell to Text *before* assigning so the
# "2026/02/28" literal isn't auto-converted into a date value, then drop
# the temporary number format again so the cell ends up unstyled, matching
# its neighbours.
$ws.Cells.Item(904, 1).NumberFormat = "@"
$ws.Cells.Item(904, 1).Value = "2026/02/28"
$ws.Cells.Item(904, 1).ClearFormats()

$ws.Cells.Item(904, 2).Value = "土"
$ws.Cells.Item(904, 3).Value = 19
$ws.Cells.Item(904, 4).Value = 201
